$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-22 04:48:16"
$ws.Range("N2").Value = "1.0 °C 4:08 TU"
$ws.Range("E3").Value = "2026-02-22 04:48:18"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "36%"
$ws.Range("K3").Value = "-0.1 MJ/m2"
$ws.Range("N3").Value = "1.4 °C 4:23 TU"
$ws.Range("O3").Value = "2.7 °C"
$ws.Range("E4").Value = "2026-02-22 04:48:20"
$ws.Range("J4").Value = "1028.2 hPa"
$ws.Range("E5").Value = "2026-02-22 04:48:23"
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "30%"
$ws.Range("N5").Value = "2.3 °C 4:15 TU"
$ws.Range("O5").Value = "4.7 °C"
$ws.Range("E6").Value = "2026-02-22 04:48:25"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "78%"
$ws.Range("N6").Value = "6.3 °C 4:04 TU"
$ws.Range("O6").Value = "7.9 °C"
$ws.Range("E7").Value = "2026-02-22 04:48:28"
$ws.Range("J7").Value = "1028.0 hPa"
$ws.Range("E8").Value = "2026-02-22 04:48:30"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "53%"
$ws.Range("J8").Value = "1027.9 hPa"
$ws.Range("M8").Value = "13.6 °C 4:28 TU"
$ws.Range("O8").Value = "11.9 °C"
$ws.Range("E9").Value = "2026-02-22 04:48:32"
$ws.Range("N9").Value = "2.5 °C 4:06 TU"
$ws.Range("O9").Value = "4.6 °C"
$ws.Range("E10").Value = "2026-02-22 04:48:33"
$ws.Range("E11").Value = "2026-02-22 04:48:34"
$ws.Range("N11").Value = "1.0 °C 4:14 TU"
$ws.Range("O11").Value = "1.6 °C"
$ws.Range("E12").Value = "2026-02-22 04:48:36"
$ws.Range("N12").Value = "2.8 °C 4:01 TU"
$ws.Range("O12").Value = "4.8 °C"
$ws.Range("E13").Value = "2026-02-22 04:48:37"
$ws.Range("O13").Value = "-2.4 °C"
$ws.Range("E14").Value = "2026-02-22 04:48:38"
$ws.Range("N14").Value = "6.8 °C 4:11 TU"
$ws.Range("O14").Value = "7.6 °C"
$ws.Range("E15").Value = "2026-02-22 04:48:39"
$ws.Range("O15").Value = "4.7 °C"
$ws.Range("E16").Value = "2026-02-22 04:48:40"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "23%"
$ws.Range("E17").Value = "2026-02-22 04:48:41"
$ws.Range("E18").Value = "2026-02-22 04:48:42"
$ws.Range("N18").Value = "0.6 °C 4:16 TU"
$ws.Range("O18").Value = "1.5 °C"
$ws.Range("E19").Value = "2026-02-22 04:48:43"
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "57%"
$ws.Range("L19").Value = "6.5 km/h - 217º 4:26 TU"
$ws.Range("E20").Value = "2026-02-22 04:48:44"
$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = "40%"
$ws.Range("O20").Value = "0.5 °C"
$ws.Range("E21").Value = "2026-02-22 04:48:46"
$ws.Range("J21").Value = "1033.5 hPa"
$ws.Range("N21").Value = "1.1 °C 4:06 TU"
$ws.Range("O21").Value = "2.8 °C"
$ws.Range("E22").Value = "2026-02-22 04:48:49"
$ws.Range("K22").Value = "-0.1 MJ/m2"
$ws.Range("L22").Value = "15.8 km/h - 319º 4:13 TU"
$ws.Range("E23").Value = "2026-02-22 04:48:51"
$ws.Range("H23").NumberFormat = "@"
$ws.Range("H23").Value = "30%"
$ws.Range("L23").Value = "9.0 km/h - 1º 4:12 TU"
$ws.Range("O23").Value = "3.7 °C"
$ws.Range("E24").Value = "2026-02-22 04:48:54"
$ws.Range("J24").Value = "1031.5 hPa"
$ws.Range("O24").Value = "1.3 °C"
$ws.Range("E25").Value = "2026-02-22 04:48:56"
$ws.Range("E26").Value = "2026-02-22 04:48:59"
$ws.Range("H26").NumberFormat = "@"
$ws.Range("H26").Value = "34%"
$ws.Range("J26").Value = "1028.6 hPa"
$ws.Range("L26").Value = "16.6 km/h - 16º 4:29 TU"
$ws.Range("E27").Value = "2026-02-22 04:49:01"
$ws.Range("H27").NumberFormat = "@"
$ws.Range("H27").Value = "31%"
$ws.Range("E28").Value = "2026-02-22 04:49:03"
$ws.Range("J28").Value = "1030.1 hPa"
$ws.Range("O28").Value = "2.3 °C"
$ws.Range("E29").Value = "2026-02-22 04:49:06"
$ws.Range("N29").Value = "2.5 °C 4:25 TU"
$ws.Range("O29").Value = "4.8 °C"
$ws.Range("E30").Value = "2026-02-22 04:49:08"
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = "82%"
$ws.Range("J30").Value = "1027.9 hPa"
$ws.Range("N30").Value = "6.7 °C 4:22 TU"
$ws.Range("O30").Value = "8.0 °C"
$ws.Range("E31").Value = "2026-02-22 04:49:11"
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = "62%"
$ws.Range("J31").Value = "1027.2 hPa"
$ws.Range("O31").Value = "12.3 °C"
$ws.Range("E32").Value = "2026-02-22 04:49:13"
$ws.Range("N32").Value = "-3.8 °C 4:29 TU"
$ws.Range("E33").Value = "2026-02-22 04:49:15"
$ws.Range("O33").Value = "1.4 °C"
$ws.Range("E34").Value = "2026-02-22 04:49:18"
$ws.Range("L34").Value = "18.4 km/h - 16º 4:00 TU"
$ws.Range("O34").Value = "2.4 °C"
$ws.Range("E35").Value = "2026-02-22 04:49:20"
$ws.Range("H35").NumberFormat = "@"
$ws.Range("H35").Value = "38%"
$ws.Range("J35").Value = "1031.7 hPa"
$ws.Range("E36").Value = "2026-02-22 04:49:23"
$ws.Range("E37").Value = "2026-02-22 04:49:25"
$ws.Range("J37").Value = "1033.9 hPa"
$ws.Range("N37").Value = "-1.2 °C 4:29 TU"
$ws.Range("O37").Value = "-0.3 °C"
$ws.Range("E38").Value = "2026-02-22 04:49:28"
$ws.Range("H38").NumberFormat = "@"
$ws.Range("H38").Value = "80%"
$ws.Range("N38").Value = "3.2 °C 4:14 TU"
$ws.Range("O38").Value = "5.3 °C"
$ws.Range("E39").Value = "2026-02-22 04:49:30"
$ws.Range("K39").Value = "-0.1 MJ/m2"
$ws.Range("E40").Value = "2026-02-22 04:49:33"
$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value = "73%"
$ws.Range("J40").Value = "1032.8 hPa"
$ws.Range("O40").Value = "4.5 °C"
$ws.Range("E41").Value = "2026-02-22 04:49:35"
$ws.Range("H41").NumberFormat = "@"
$ws.Range("H41").Value = "92%"
$ws.Range("J41").Value = "1028.6 hPa"
$ws.Range("K41").Value = "-0.1 MJ/m2"
$ws.Range("N41").Value = "4.0 °C 4:12 TU"
$ws.Range("O41").Value = "5.4 °C"
$ws.Range("E42").Value = "2026-02-22 04:49:37"
$ws.Range("N42").Value = "3.6 °C 4:29 TU"
$ws.Range("O42").Value = "5.1 °C"
$ws.Range("E43").Value = "2026-02-22 04:49:39"
$ws.Range("N43").Value = "0.8 °C 4:11 TU"
$ws.Range("O43").Value = "2.3 °C"
$ws.Range("E44").Value = "2026-02-22 04:49:42"
$ws.Range("N44").Value = "-1.3 °C 4:19 TU"
$ws.Range("O44").Value = "0.1 °C"
$ws.Range("E45").Value = "2026-02-22 04:49:44"
$ws.Range("H45").NumberFormat = "@"
$ws.Range("H45").Value = "61%"
$ws.Range("J45").Value = "1031.6 hPa"
$ws.Range("N45").Value = "2.8 °C 4:26 TU"
$ws.Range("O45").Value = "4.7 °C"
$ws.Range("E46").Value = "2026-02-22 04:49:46"
$ws.Range("J46").Value = "1031.3 hPa"
$ws.Range("N46").Value = "0.0 °C 4:29 TU"
$ws.Range("O46").Value = "1.4 °C"
